# Update cryptocurrency price/volume data as of the latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.663.94"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.588.90"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.23%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "208.43"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  +0.27%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "22.26"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("E9").Value = "  -1.90%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.0592"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0867"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "1.813.83"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "1.577.71"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "27.638.44"
$ws.Range("E16").Value = "  -0.97%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "63.41"
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "219.95"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -3.54%  "
$ws.Range("D19").Value = "0.0₃0697"
$ws.Range("E19").Value = "  -3.00%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.33"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("E21").Value = "  +0.32%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "4.14"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -4.55%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "9.69"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("E24").Value = "  -3.32%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "153.66"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "6.83"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  +0.26%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "15.13"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  -4.70%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.16"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.0469"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").Value = "1.377.70"
$ws.Range("E33").Value = "  -2.62%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.94"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("E35").Value = "  -4.76%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.976"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("E37").Value = "  -0.68%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.0168"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.536"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("E41").Value = "  +0.32%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.970"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -3.36%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "64.27"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("E44").Value = "  +2.12%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "5.24"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "1.724.91"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("E47").Value = "  -5.17%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "87.23"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -0.94%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.0966"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -4.20%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0494"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -1.58%  "
